$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Miscellaneous" section (column C/D starting row 18) gains Midnight Commander
# entries: fte, mc, mcedit (plus their descriptions), so it's available on all
# releases.
$ws.Range("C23").Value = "fte"
$ws.Range("D23").Value = "IDE/Text Editor"

$ws.Range("C21").Value = "mc"
$ws.Range("C22").Value = "mcedit"
$ws.Range("D22").Value = "Text Editor"
$ws.Range("D21").Value = "File Explorer (Norton Commander Like)"

$ws.Rows("22:23").AutoFit()

$ws.Range("C24:D24").Clear()

$ws.Range("D19").Select()
